# Applies the OOXML diff to the attendance/session-analysis workbook.
#
# Summary of the change (per the commit "Sync attendance_reports, ... - 2026-01-04"):
#  1. The "Recorded By" text in column G that lists both the instructor and
#     "System" has the two names swapped everywhere:
#        "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#  2. The reporting "today" moved forward to 04/01/2026, so the session that
#     was dated 04/01/2026 for each of the six B1 sub-groups (rows 20, 46,
#     72, 98, 124, 150) is no longer "Pending" but "Not Recorded" - its
#     status text and its row fill/style change accordingly (matching the
#     existing "Not Recorded" look used e.g. in row 3).
#  3. Because six sessions moved from the "Pending" bucket into the
#     "Missing" bucket, the overall summary counters update:
#        L7 (Missing Sessions): 9   -> 15
#        L8 (Pending Sessions): 102 -> 96
#     and, per affected group (rows 15-20 of the "Group Statistics" table),
#     the Missing (P) column goes up by 1 and the Pending (Q) column goes
#     down by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#    everywhere it occurs in column G (the "Recorded By" column).
# ---------------------------------------------------------------------
$recordedByRange = $ws.Range("G1:G319")
$recordedByRange.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Flip the six "04/01/2026" session rows from Pending -> Not Recorded.
# ---------------------------------------------------------------------
$pendingRows = @(20, 46, 72, 98, 124, 150)

foreach ($r in $pendingRows) {
    # Update the status text shown in column I.
    $ws.Cells.Item($r, 9).Value = "Not Recorded"

    # Re-style columns A:I to match the existing "Not Recorded" look
    # (copy the formatting already used by row 3, which has that status).
    $srcFormat = $ws.Range("A3:I3")
    $dstFormat = $ws.Range("A" + $r + ":I" + $r)
    $srcFormat.Copy() | Out-Null
    $dstFormat.PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Update the summary counters.
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 15
$ws.Range("L8").Value = 96

# Per-group Missing (P) / Pending (Q) counts for rows 15-20
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 7

$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 7

$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 7

$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 7

$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 7

$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 7
